# This workbook is an ObjTables/SBtab export. Every worksheet carries one (or
# two, for the very first sheet) header row(s) of the form:
#   !!!ObjTables schema='SBtab' objTablesVersion='<ver>' date='<date>'
#   !!ObjTables schema='SBtab' type='Data' tableFormat='row' class='<Class>' name='<Class>' date='<date>' objTablesVersion='<ver>'
#
# The commit bumps the generation timestamp/version embedded in those header
# strings from date='2020-04-27 01:07:46' / objTablesVersion='0.0.9' to
# date='2020-05-29 00:21:41' / objTablesVersion='1.0.0' everywhere they occur.

$wb = $excel.ActiveWorkbook

$oldDate = "2020-04-27 01:07:46"
$newDate = "2020-05-29 00:21:41"
$oldVer  = "0.0.9"
$newVer  = "1.0.0"

# The sheets are protected (no password), which blocks writing to cells.
# Unprotect every sheet first so the header cells can be updated.
foreach ($ws in $wb.Worksheets) {
    $ws.Unprotect()
}

foreach ($ws in $wb.Worksheets) {
    # Row 1 on every sheet holds that sheet's own "!!ObjTables ... class='X'"
    # header (the very first sheet's A1 instead holds the workbook-level
    # "!!!ObjTables ..." banner) - update whichever text is present.
    $cellA1 = $ws.Range("A1")
    $textA1 = $cellA1.Text
    if ($textA1 -like "*ObjTables*date=*") {
        $updated = $textA1 -replace [regex]::Escape($oldDate), $newDate
        $updated = $updated -replace [regex]::Escape($oldVer), $newVer
        $cellA1.Value = $updated
    }

    # Only the first sheet (!!Compartment) additionally has the per-class
    # header duplicated into A2; update it too if present.
    $cellA2 = $ws.Range("A2")
    $textA2 = $cellA2.Text
    if ($textA2 -like "*ObjTables*date=*") {
        $updated2 = $textA2 -replace [regex]::Escape($oldDate), $newDate
        $updated2 = $updated2 -replace [regex]::Escape($oldVer), $newVer
        $cellA2.Value = $updated2
    }
}

# Restore sheet protection on every sheet.
foreach ($ws in $wb.Worksheets) {
    $ws.Protect([System.Reflection.Missing]::Value, $true, $true, $true)
}
